$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.252669095993042
$ws.Range("B1").Value = 2.487767457962036
$ws.Range("C1").Value = 3.509603261947632
$ws.Range("D1").Value = 3.651207208633423
$ws.Range("E1").Value = 0.8329172730445862
